# Update "want to go" counts (column F) across the three sheets that contain
# this data: 展览 (Exhibition), 演出 (Performance), 全部类型 (All Types).
# 本地生活 has no data rows, so it is untouched.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 5214
$ws.Range("F4").Value  = 5214
$ws.Range("F6").Value  = 166
$ws.Range("F7").Value  = 214
$ws.Range("F10").Value = 179
$ws.Range("F11").Value = 8706
$ws.Range("F12").Value = 277
$ws.Range("F15").Value = 9
$ws.Range("F16").Value = 2585
$ws.Range("F20").Value = 9
$ws.Range("F24").Value = 6524
$ws.Range("F25").Value = 210
$ws.Range("F27").Value = 147
$ws.Range("F30").Value = 7088
$ws.Range("F33").Value = 234
$ws.Range("F38").Value = 7
$ws.Range("F41").Value = 2545
$ws.Range("F47").Value = 546
$ws.Range("F48").Value = 3507
$ws.Range("F49").Value = 96
$ws.Range("F50").Value = 1131

# --- 演出 (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 15
$ws.Range("F3").Value = 2
$ws.Range("F5").Value = 198
$ws.Range("F7").Value = 87

# --- 全部类型 (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 5214
$ws.Range("F4").Value  = 5214
$ws.Range("F6").Value  = 166
$ws.Range("F7").Value  = 214
$ws.Range("F9").Value  = 179
$ws.Range("F10").Value = 8706
$ws.Range("F11").Value = 8706
$ws.Range("F12").Value = 277
$ws.Range("F15").Value = 9
$ws.Range("F16").Value = 15
$ws.Range("F17").Value = 2585
$ws.Range("F18").Value = 2
$ws.Range("F20").Value = 198
$ws.Range("F23").Value = 87
$ws.Range("F24").Value = 9
$ws.Range("F29").Value = 6524
$ws.Range("F30").Value = 210
$ws.Range("F32").Value = 147
$ws.Range("F34").Value = 7088
$ws.Range("F38").Value = 7
$ws.Range("F45").Value = 546
$ws.Range("F47").Value = 3507
$ws.Range("F48").Value = 96
$ws.Range("F50").Value = 1131
